$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove stray y_0_forecast values for the first two rows (naive forecaster
# bug: these shouldn't have been populated).
$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()

# Recomputed forecast values after the bug fix (tiny precision corrections).
$ws.Range("C4").Value = -14.45332333832744
$ws.Range("C5").Value = 8.600536527919612
$ws.Range("E5").Value = 9.002271992040356
$ws.Range("E6").Value = 7.550992341868912
$ws.Range("C7").Value = 4.639893381363192
$ws.Range("E7").Value = 5.799303245920906
$ws.Range("E8").Value = 2.638010271840918
$ws.Range("E11").Value = 2.129835064860441
$ws.Range("C12").Value = 4.695933104194361
$ws.Range("C17").Value = 5.120680133083622
$ws.Range("E17").Value = 3.315588342229514
$ws.Range("C18").Value = -0.5532735011319123
